$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2190.35
$ws.Range("I40").Value = 2933.111
$ws.Range("J40").Value = 1974.7097
$ws.Range("K40").Value = 2933.111
$ws.Range("L40").Value = 1974.7097
$ws.Range("M40").Value = -2758.111
$ws.Range("N40").Value = -2324.7097
$ws.Range("H80").Value = 13944791
$ws.Range("I80").Value = 25020166
$ws.Range("J80").Value = 100573.25
$ws.Range("K80").Value = 75060498
$ws.Range("L80").Value = 301719.75
$ws.Range("M80").Value = -75059500
$ws.Range("N80").Value = -303715.75
$ws.Range("H83").Value = 13944791
$ws.Range("I83").Value = 25020166
$ws.Range("J83").Value = 100573.25
$ws.Range("K83").Value = 225181494
$ws.Range("L83").Value = 905159.25
$ws.Range("M83").Value = -225176502
$ws.Range("N83").Value = -915143.25
$ws.Range("H88").Value = 11935449
$ws.Range("I88").Value = 33336142
$ws.Range("K88").Value = 33336142
$ws.Range("M88").Value = -33335736
$ws.Range("H91").Value = 11935449
$ws.Range("I91").Value = 33336142
$ws.Range("K91").Value = 33336142
$ws.Range("M91").Value = -33334738

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2237328
$ws.Range("I32").Value = 2319126.8
$ws.Range("K32").Value = 2319126.8
$ws.Range("M32").Value = -2318839.8
$ws.Range("H45").Value = 5413
$ws.Range("I45").Value = 1797.3334
$ws.Range("J45").Value = 7220.8335
$ws.Range("K45").Value = 1797.3334
$ws.Range("L45").Value = 7220.8335
$ws.Range("M45").Value = -1420.3334
$ws.Range("N45").Value = -7974.8335
$ws.Range("H110").Value = 16667645
$ws.Range("I110").Value = 819.5333000000001
$ws.Range("K110").Value = 819.5333000000001
$ws.Range("M110").Value = 1225.4667
$ws.Range("H132").Value = 3545.56
$ws.Range("I132").Value = 1398.0834
$ws.Range("K132").Value = 4194.2502
$ws.Range("M132").Value = -1664.2502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 33786764
$ws.Range("I86").Value = 13891377
$ws.Range("K86").Value = 13891377
$ws.Range("M86").Value = -13890254
$ws.Range("H89").Value = 33786764
$ws.Range("I89").Value = 13891377
$ws.Range("K89").Value = 69456885
$ws.Range("M89").Value = -69451269
$ws.Range("H107").Value = 46878468
$ws.Range("I107").Value = 59212924
$ws.Range("J107").Value = 7522.4
$ws.Range("K107").Value = 59212924
$ws.Range("L107").Value = 7522.4
$ws.Range("M107").Value = -59211004
$ws.Range("N107").Value = -11362.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11150.219
$ws.Range("I31").Value = 5299
$ws.Range("J31").Value = 14660.95
$ws.Range("K31").Value = 5299
$ws.Range("L31").Value = 14660.95
$ws.Range("M31").Value = -5004
$ws.Range("N31").Value = -15250.95
$ws.Range("H34").Value = 11150.219
$ws.Range("I34").Value = 5299
$ws.Range("J34").Value = 14660.95
$ws.Range("K34").Value = 5299
$ws.Range("L34").Value = 14660.95
$ws.Range("M34").Value = -5097
$ws.Range("N34").Value = -15064.95
$ws.Range("H58").Value = 17249370
$ws.Range("I58").Value = 62500708
$ws.Range("K58").Value = 62500708
$ws.Range("M58").Value = -62500505
$ws.Range("H107").Value = 1948.375
$ws.Range("I107").Value = 1667.5555
$ws.Range("J107").Value = 2309.4285
$ws.Range("K107").Value = 1667.5555
$ws.Range("L107").Value = 2309.4285
$ws.Range("M107").Value = 252.4445000000001
$ws.Range("N107").Value = -6149.4285
$ws.Range("H134").Value = 7004.433
$ws.Range("I134").Value = 1371
$ws.Range("K134").Value = 4113
$ws.Range("M134").Value = -1578
$ws.Range("H136").Value = 17249370
$ws.Range("I136").Value = 62500708
$ws.Range("K136").Value = 187502124
$ws.Range("M136").Value = -187499574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 23809776
$ws.Range("I33").Value = 41666860
$ws.Range("K33").Value = 250001160
$ws.Range("M33").Value = -250000877
$ws.Range("H92").Value = 8548563
$ws.Range("J92").Value = 8548563
$ws.Range("L92").Value = 25645689
$ws.Range("N92").Value = -25648185
$ws.Range("H113").Value = 4729.25
$ws.Range("J113").Value = 6500.8
$ws.Range("L113").Value = 19502.4
$ws.Range("N113").Value = -23842.4
$ws.Range("H131").Value = 1746.3077
$ws.Range("J131").Value = 3338.6
$ws.Range("L131").Value = 10015.8
$ws.Range("N131").Value = -20095.8
$ws.Range("H139").Value = 66458.766
$ws.Range("I139").Value = 254324.75
$ws.Range("J139").Value = 8653.846
$ws.Range("K139").Value = 762974.25
$ws.Range("L139").Value = 25961.538
$ws.Range("M139").Value = -757834.25
$ws.Range("N139").Value = -36241.538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3554.5854
$ws.Range("I132").Value = 1405.9231
$ws.Range("J132").Value = 7278.933
$ws.Range("K132").Value = 4217.7693
$ws.Range("L132").Value = 21836.799
$ws.Range("M132").Value = -1687.7693
$ws.Range("N132").Value = -26896.799

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4628.857
$ws.Range("I68").Value = 1250.5
$ws.Range("J68").Value = 9133.333000000001
$ws.Range("K68").Value = 1250.5
$ws.Range("L68").Value = 9133.333000000001
$ws.Range("M68").Value = -501.5
$ws.Range("N68").Value = -10631.333
$ws.Range("H71").Value = 4628.857
$ws.Range("I71").Value = 1250.5
$ws.Range("J71").Value = 9133.333000000001
$ws.Range("K71").Value = 6252.5
$ws.Range("L71").Value = 45666.665
$ws.Range("M71").Value = -2508.5
$ws.Range("N71").Value = -53154.665
$ws.Range("H122").Value = 6754.154
$ws.Range("I122").Value = 3766.6667
$ws.Range("J122").Value = 7650.4
$ws.Range("K122").Value = 11300.0001
$ws.Range("L122").Value = 22951.2
$ws.Range("M122").Value = -8850.000100000001
$ws.Range("N122").Value = -27851.2
$ws.Range("H132").Value = 13896822
$ws.Range("I132").Value = 41669984
$ws.Range("K132").Value = 125009952
$ws.Range("M132").Value = -125007422

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8080957
$ws.Range("J81").Value = 12505330
$ws.Range("L81").Value = 25010660
$ws.Range("N81").Value = -25012782
$ws.Range("H84").Value = 8080957
$ws.Range("J84").Value = 12505330
$ws.Range("L84").Value = 125053300
$ws.Range("N84").Value = -125063908
$ws.Range("H136").Value = 41674020
$ws.Range("J136").Value = 9050.166999999999
$ws.Range("L136").Value = 27150.501
$ws.Range("N136").Value = -32250.501
